$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
  0.26837214234553386,
  -0.06229747168520916,
  -0.0039999999674176223,
  -0.0079999999419353429,
  -0.0029999999654295451,
  -0.0019999999613578012,
  -0.0099999999167765807,
  -0.0099999999145690133,
  -0.0019999999556157277,
  -0.0019999999533446555,
  0.0061274791962206976,
  -0.0034999999441360785,
  -0.0034999999405922466,
  -0.0079999999150484058,
  -0.00099999995255029006,
  -0.0019999999464803686,
  -0.0019999999455597717,
  -0.0039999999343862669,
  -0.0039999999761510807,
  -0.0039999999744484427,
  -0.0039999999741731074,
  -0.0039999999739439573,
  -0.031666237181782186,
  -0.019999999874287688,
  -0.019999999872634788,
  -0.018635793561399439,
  -0.0024999999570654019,
  -0.0019999999561717274,
  -0.0069999999264753754,
  -0.059999999637014856,
  0.050438407403776253,
  -0.0099999999120434779,
  -0.0039999999446038714
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
